# Fixed up incorrect import code, closes #232
#
# This script fixes the sample data in the "Assignments" and
# "Participant Experiments" sheets of the import template so that the
# participant_experiment:assignments values actually reference assignment
# ids that exist in the Assignments sheet, and adds the extra
# assignment rows those references need.

$wb = $excel.ActiveWorkbook
$wsAssignments = $wb.Worksheets.Item(2)
$wsParticipantExperiments = $wb.Worksheets.Item(3)

# --- Assignments sheet -------------------------------------------------
# Correct the existing example rows (2 and 3 in the media_items/activity/id
# table) and append three brand-new rows so ids 4, 5 and 6 exist for the
# Participant Experiments sheet to reference.

$wsAssignments.Range("A3").Value = 3
$wsAssignments.Range("A4").Value = 2

# Copy the number formatting from the last existing data row down onto the
# new rows/cells before filling in their values.
$wsAssignments.Range("A4:C4").Copy()
$wsAssignments.Range("A5:C7").PasteSpecial(-4122)

$wsAssignments.Range("C4").Copy()
$wsAssignments.Range("C8").PasteSpecial(-4122)
$wsAssignments.Range("C9").PasteSpecial(-4122)

$wsAssignments.Range("A5").Value = 3
$wsAssignments.Range("B5").Value = 1
$wsAssignments.Range("C5").Value = 4

$wsAssignments.Range("A6").Value = 2
$wsAssignments.Range("B6").Value = 2
$wsAssignments.Range("C6").Value = 5

$wsAssignments.Range("A7").Value = 1
$wsAssignments.Range("B7").Value = 3
$wsAssignments.Range("C7").Value = 6

# --- Participant Experiments sheet --------------------------------------
# Widen the first column a bit and fix up the assignment id lists so they
# point at real assignment ids (1-6) instead of the previous bogus
# combinations.

$wsParticipantExperiments.Columns.Item(1).ColumnWidth = 23.33

$wsParticipantExperiments.Range("B3").Value = "3,4"
$wsParticipantExperiments.Range("B4").Value = "5,6"
